# "Database structure.xlsx" — split the single "Reporting module" row on the
# Modules sheet into three separate rows (School level / Block level /
# Education Department level), renumber the Module ID column, and update the
# dependent defined name / used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Modules")

# Row 18 currently holds Module ID 17 "Reporting module" (Phase 1, Yes).
# Insert two blank rows right after it so it becomes three rows in total.
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(19).Insert()

# Row 18: rename the original "Reporting module" entry.
$ws.Cells.Item(18, 2).Value = "Reporting module - School level"
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 4).Value = "Yes"

# Row 19: new "Block level" entry.
$ws.Cells.Item(19, 2).Value = "Reporting module - Block level"
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = "No"

# Row 20: new "Education Department level" entry.
$ws.Cells.Item(20, 2).Value = "Reporting module - Education Department level"
$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(20, 4).Value = "No"

# Renumber the Module ID column (A) for every data row so it stays
# sequential (1..30) after the split.
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Keep the hidden AutoFilter defined name in sync with the new data extent
# (header + 28 rows -> up to row 29 of the range, i.e. $C$29 as in the diff).
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Modules!`$A`$1:`$C`$29"
